# Automatische test-sync: 2025-08-02 00:19:50
# Adds the newest "Testmail #20" log row to the Logs sheet, rolls the new
# "Klacht / Probleem" category into the Dashboard pivot-style summary table,
# and extends the conditional formatting + chart series ranges so the new
# rows are covered.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Logs sheet: append row 10 with the new test-mail entry
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A10").Value = "Ik ben niet tevreden over hoe dit is gegaan."
$logs.Range("B10").Value = "mailmind.test@zohomail.eu"
$logs.Range("C10").Value = "Testmail #20: Ik ben niet tevreden over hoe dit is gegaan."
$logs.Range("D10").Value = "Klacht / Probleem"
$logs.Range("E10").Value = "Bedankt, we hebben dit doorgestuurd naar support@bedrijf.nl."
$logs.Range("F10").Value = "2025-08-02 00:18:52"
$logs.Range("G10").Value = "Ja"
$logs.Range("H10").Value = "Ja"
$logs.Range("I10").Value = "Nee"
$logs.Range("J10").Value = "Nee"

# Conditional formatting on the Logs sheet is keyed to the old 2:9 row
# range; grow every rule's AppliesTo range down to the new row 10 so the
# whole block (all cfRules sharing the sqref) moves together.
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "9")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "10")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count(); $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 2) Dashboard sheet: roll the new category into the summary table
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A7").Value = "Klacht / Probleem"
$dash.Range("B7").Value = 1

# ---------------------------------------------------------------------
# 3) Chart on the Dashboard sheet: extend category/value series ranges
# ---------------------------------------------------------------------
$chart = $dash.ChartObjects().Item(1).Chart()
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$7,'Dashboard'!`$B`$2:`$B`$7,1)"
